$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '95.957.73'
$ws.Range("E2").Value = '  +1.66%  '

# Row 3
$ws.Range("D3").Value = '3.594.42'
$ws.Range("E3").Value = '  +3.29%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.09'
$ws.Range("E5").Value = '  +0.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '655.01'
$ws.Range("E6").Value = '  +4.73%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.47'
$ws.Range("E7").Value = '  +1.84%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.403'
$ws.Range("E8").Value = '  +2.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.02%  '

# Row 10
$ws.Range("E10").Value = '  -0.47%  '

# Row 11
$ws.Range("D11").Value = '3.592.16'
$ws.Range("E11").Value = '  +3.27%  '

# Row 12
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.73'
$ws.Range("E12").Value = '  -1.39%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.202'
$ws.Range("E13").Value = '  +0.61%  '

# Row 14
$ws.Range("E14").Value = '  +3.73%  '

# Row 15
$ws.Range("D15").Value = '4.291.64'
$ws.Range("E15").Value = '  +3.67%  '

# Row 16
$ws.Range("D16").Value = '95.820.52'
$ws.Range("E16").Value = '  +1.71%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000254'
$ws.Range("E17").Value = '  +1.63%  '

# Row 18
$ws.Range("D18").Value = '3.599.23'
$ws.Range("E18").Value = '  +3.48%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.94'
$ws.Range("E19").Value = '  -5.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.72'
$ws.Range("E20").Value = '  +0.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.96'
$ws.Range("E21").Value = '  -0.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.54'
$ws.Range("E22").Value = '  +5.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '509.94'
$ws.Range("E23").Value = '  -2.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.481'
$ws.Range("E24").Value = '  -5.38%  '

# Row 25
$ws.Range("E25").Value = '  +6.05%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.37'
$ws.Range("E26").Value = '  -5.17%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '92.34'
$ws.Range("E27").Value = '  -4.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.60'
$ws.Range("E28").Value = '  +3.35%  '

# Row 29
$ws.Range("D29").Value = '3.793.03'
$ws.Range("E29").Value = '  +3.49%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.09'
$ws.Range("E30").Value = '  +5.56%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.31'
$ws.Range("E31").Value = '  -1.43%  '

# Row 32
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.07%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.140'
$ws.Range("E33").Value = '  -1.25%  '

# Row 34
$ws.Range("E34").Value = '  +1.65%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.42'
$ws.Range("E35").Value = '  +9.20%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.179'
$ws.Range("E36").Value = '  -0.56%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.565'
$ws.Range("E37").Value = '  +0.72%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.18'
$ws.Range("E38").Value = '  +8.21%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '561.44'
$ws.Range("E39").Value = '  -2.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.48'
$ws.Range("E40").Value = '  +2.45%  '

# Row 41
$ws.Range("E41").Value = '  +0.02%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.151'
$ws.Range("E42").Value = '  +0.88%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.911'
$ws.Range("E43").Value = '  -1.68%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '35.55'
$ws.Range("E44").Value = '  +40.47%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.76'
$ws.Range("E45").Value = '  +3.41%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.31'
$ws.Range("E46").Value = '  +6.98%  '

# Row 47
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.71'
$ws.Range("E47").Value = '  +3.30%  '

# Row 48
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.60'
$ws.Range("E48").Value = '  -0.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0414'
$ws.Range("E49").Value = '  -2.66%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.58'
$ws.Range("E50").Value = '  +0.92%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.78'
$ws.Range("E51").Value = '  +0.98%  '
